{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Implements the three user-visible content changes from the diff:\n//   1. \"Dashboard Automation project is\" -> \"DM01 Data Processing Project is\"\n//   2. Removes \" and create an output in the form of a Dashboard and a\n//      PowerPoint Presentation\" (the sentence now ends after \"reporting\n//      pipeline.\")\n//   3. Turns the lone \"--\" paragraph (under \"Clone the GitHub Repository\")\n//      into a hyperlink pointing at the DM01_data_processing GitHub repo,\n//      followed by a trailing space run.\n\nconst body = context.document.body;\n\n// --- 1. Rename the project in the intro paragraph -------------------------\nconst titleHits = body.search(\"Dashboard Automation project is\", { matchCase: true });\ntitleHits.load(\"text\");\nawait context.sync();\n\nif (titleHits.items.length > 0) {\n  titleHits.items[0].insertText(\"DM01 Data Processing Project is\", \"Replace\");\n  await context.sync();\n}\n\n// --- 2. Drop the \"...and create an output...\" clause ----------------------\nconst clauseHits = body.search(\n  \" and create an output in the form of a Dashboard and a PowerPoint Presentation\",\n  { matchCase: true }\n);\nclauseHits.load(\"text\");\nawait context.sync();\n\nif (clauseHits.items.length > 0) {\n  clauseHits.items[0].insertText(\"\", \"Replace\");\n  await context.sync();\n}\n\n// --- 3. Turn the placeholder \"--\" into the GitHub repo hyperlink ----------\nconst repoUrl = \"https://github.com/Arush313/DM01_data_processing/tree/main\";\n\nconst placeholderHits = body.search(\"--\", { matchCase: true });\nplaceholderHits.load(\"text\");\nawait context.sync();\n\nif (placeholderHits.items.length > 0) {\n  // Replace the placeholder text (plain) with the URL + a trailing space,\n  // then re-find just the URL portion so the hyperlink formatting/relationship\n  // only wraps the link text, leaving the trailing space as a plain run\n  // (matches the target markup: <w:hyperlink>\u2026</w:hyperlink><w:r> </w:r>).\n  placeholderHits.items[0].insertText(repoUrl + \" \", \"Replace\");\n  await context.sync();\n\n  const urlHits = body.search(repoUrl, { matchCase: true });\n  urlHits.load(\"text\");\n  await context.sync();\n\n  if (urlHits.items.length > 0) {\n    urlHits.items[0].hyperlink = repoUrl;\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Implements the three user-visible content changes from the diff:\n#   1. \"Dashboard Automation project is\" -> \"DM01 Data Processing Project is\"\n#   2. Removes \" and create an output in the form of a Dashboard and a\n#      PowerPoint Presentation\" (the sentence now ends after \"reporting\n#      pipeline.\")\n#   3. Turns the lone \"--\" paragraph (under \"Clone the GitHub Repository\")\n#      into a hyperlink pointing at the DM01_data_processing GitHub repo,\n#      followed by a trailing space run.\n\n$d = $word.ActiveDocument\n\n# --- 1. Rename the project in the intro paragraph --------------------------\n$find1 = $d.Content.Find\n$find1.Execute(\n    \"Dashboard Automation project is\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"DM01 Data Processing Project is\", 2\n) | Out-Null\n\n# --- 2. Drop the \"...and create an output...\" clause -----------------------\n$find2 = $d.Content.Find\n$find2.Execute(\n    \" and create an output in the form of a Dashboard and a PowerPoint Presentation\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"\", 2\n) | Out-Null\n\n# --- 3. Turn the placeholder \"--\" into the GitHub repo hyperlink -----------\n$repoUrl = \"https://github.com/Arush313/DM01_data_processing/tree/main\"\n\n$linkRange = $d.Content\n$find3 = $linkRange.Find\n$found = $find3.Execute(\n    \"--\",\n    $true, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"\", 0\n)\n\nif ($found) {\n    # $linkRange now spans exactly the \"--\" placeholder text (Find.Execute\n    # collapses/extends Range to the hit). Replace it with the URL text and\n    # wrap that same range with a hyperlink relationship.\n    $linkRange.Text = $repoUrl\n    $d.Hyperlinks.Add($linkRange, $repoUrl) | Out-Null\n\n    # Insert a trailing plain-text space right after the hyperlink, matching\n    # <w:hyperlink>...</w:hyperlink><w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n    $afterRange = $d.Range($linkRange.End, $linkRange.End)\n    $afterRange.InsertAfter(\" \")\n}\n"}
